$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.568.27'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.992.60'
$ws.Range('E3').Value = '  +6.12%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.82'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4684'
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3954'
$ws.Range('E8').Value = '  +2.15%  '
$ws.Range('E9').Value = '  -0.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08138'
$ws.Range('E10').Value = '  +3.73%  '
$ws.Range('E11').Value = '  +1.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.93'
$ws.Range('E12').Value = '  +5.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.994.69'
$ws.Range('E13').Value = '  +3.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.246'
$ws.Range('E14').Value = '  +3.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.873'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07133'
$ws.Range('E16').Value = '  +2.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.84'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001005'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.45'
$ws.Range('E20').Value = '  +3.14%  '
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.578.87'
$ws.Range('E22').Value = '  +2.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.560'
$ws.Range('E24').Value = '  +2.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.114'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '157.86'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.66'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.983'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '120.42'
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.945'
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09460'
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9160'
$ws.Range('E32').Value = '  +1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.282'
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.354'
$ws.Range('E34').Value = '  +2.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.184'
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05859'
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.180'
$ws.Range('E37').Value = '  +0.98%  '
$ws.Range('E38').Value = '  +2.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.000003370'
$ws.Range('E39').Value = '  +77.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.906'
$ws.Range('E40').Value = '  +3.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5791'
$ws.Range('E41').Value = '  +2.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1826'
$ws.Range('E42').Value = '  +3.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.894'
$ws.Range('E43').Value = '  +2.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.811'
$ws.Range('E44').Value = '  +10.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.06'
$ws.Range('E45').Value = '  +0.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5401'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06968'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.872'
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '114.11'
$ws.Range('E50').Value = '  +1.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3090'
$ws.Range('E51').Value = '  +8.18%  '
